# Update NATMI Lama1-Itga1 LR-pair TPM-derived statistics.
# The underlying TPM values for the receptor (Itga1) in the MuSCs target
# cluster were recomputed, which changes:
#   - Receptor average/total expression value for MuSCs (M/N on rows 2 & 5)
#   - Receptor derived specificity of average/total expression value,
#     recalculated across all three target clusters (O/P on rows 2-7)
#   - Edge average/total expression weight (Q/R), which depends on the
#     receptor average/total expression value (rows 2 & 5 target MuSCs)
#   - Edge average/total expression derived specificity (S/T on rows 2-7)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending=FAPs, Target=MuSCs)
$ws.Range("M2").Value = 72.07569866666667
$ws.Range("N2").Value = 216.227096
$ws.Range("O2").Value = 0.4479522040449755
$ws.Range("P2").Value = 0.4479522040449755
$ws.Range("Q2").Value = 1.423278821570667
$ws.Range("R2").Value = 12.809509394136
$ws.Range("S2").Value = 0.3199172576230066
$ws.Range("T2").Value = 0.3199172576230066

# Row 3 (Sending=FAPs, Target=FAPs)
$ws.Range("O3").Value = 0.04737448730867841
$ws.Range("P3").Value = 0.0473744873086784
$ws.Range("S3").Value = 0.03383377941715995
$ws.Range("T3").Value = 0.03383377941715994

# Row 4 (Sending=FAPs, Target=ECs)
$ws.Range("O4").Value = 0.5046733086463462
$ws.Range("P4").Value = 0.5046733086463461
$ws.Range("Q4").Value = 1.603498823138334
$ws.Range("S4").Value = 0.360426178149707
$ws.Range("T4").Value = 0.3604261781497069

# Row 5 (Sending=ECs, Target=MuSCs)
$ws.Range("M5").Value = 72.07569866666667
$ws.Range("N5").Value = 216.227096
$ws.Range("O5").Value = 0.4479522040449755
$ws.Range("P5").Value = 0.4479522040449755
$ws.Range("Q5").Value = 0.5696142465626667
$ws.Range("R5").Value = 5.126528219064
$ws.Range("S5").Value = 0.1280349464219689
$ws.Range("T5").Value = 0.128034946421969

# Row 6 (Sending=ECs, Target=FAPs)
$ws.Range("O6").Value = 0.04737448730867841
$ws.Range("P6").Value = 0.0473744873086784
$ws.Range("S6").Value = 0.01354070789151846
$ws.Range("T6").Value = 0.01354070789151846

# Row 7 (Sending=ECs, Target=ECs)
$ws.Range("O7").Value = 0.5046733086463462
$ws.Range("P7").Value = 0.5046733086463461
$ws.Range("S7").Value = 0.1442471304966392
$ws.Range("T7").Value = 0.1442471304966392
